$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch (empty) row 3 so it keeps spanning the full header width, the
# same way it did before the new columns were added. Explicitly
# (re)setting a border to "none" registers the cell in the saved XML
# without altering its (default) appearance or creating a new style.
$ws.Range("A3").Borders.Item(7).LineStyle = -4142

# ------------------------------------------------------------------
# 1) Insert 3 new columns between "Colonia *" (R) and "Calle *" (old S)
#    -> Cve Colonia, Cve Interventor, Cve Tipo Calle
#    New columns inherit the format of column R (style index 2).
# ------------------------------------------------------------------
$ws.Range("S1:U1").EntireColumn.Insert()
$ws.Range("S4").Value = "Cve Colonia"
$ws.Range("T4").Value = "Cve Interventor"
$ws.Range("U4").Value = "Cve Tipo Calle"

# ------------------------------------------------------------------
# 2) Insert 6 new columns right before "RESPONSABLE DE VALIDACION"
#    (which, after step 1, lives in column AU)
#    -> ENLACE INTERVENCION 1/2/3, FECHA SOLICITUD,
#       RESPONSABLE DE LA ENTREGA, ESTATUS ORIGEN
# ------------------------------------------------------------------
$ws.Range("AU1:AZ1").EntireColumn.Insert()

# The 6 new header cells should use the same look as the rest of the
# plain header cells (style index 2, same as column R/S/T/U above),
# not the style inherited from their left neighbour on insert.
$ws.Range("R4").Copy()
$ws.Range("AU4:AZ4").PasteSpecial(-4122)

$ws.Range("AU4").Value = "ENLACE INTERVENCION 1"
$ws.Range("AV4").Value = "ENLACE INTERVENCION 2"
$ws.Range("AW4").Value = "ENLACE INTERVENCION 3"
$ws.Range("AX4").Value = "FECHA SOLICITUD"
$ws.Range("AY4").Value = "RESPONSABLE DE LA ENTREGA"
$ws.Range("AZ4").Value = "ESTATUS ORIGEN"

# ------------------------------------------------------------------
# 3) "RESPONSABLE DE VALIDACION" / "FECHA VALIDACION" (now in BA / BB)
#    get a new highlight fill (theme accent6, darker 25%, white pattern
#    color) while keeping the same font/border/alignment as before.
# ------------------------------------------------------------------
$ws.Range("R4").Copy()
$ws.Range("BA4:BB4").PasteSpecial(-4122)
$ws.Range("BA4:BB4").Interior.ThemeColor = 10
$ws.Range("BA4:BB4").Interior.TintAndShade = -0.249977111117893
$ws.Range("BA4:BB4").Interior.PatternColor = 16777215

$excel.CutCopyMode = 0

# Leave the same cell selected/active as in the authored workbook.
$ws.Range("BA9").Select() | Out-Null
